# Add the "SpikeList Output Variable" enum entry (SPIKE_LIST_REQ) as a new
# row inserted immediately before the old "INITIAL_STATE_REQ" row, keeping
# its exponent value (14) and pushing INITIAL_STATE_REQ/FINAL_STATE_REQ
# down one row each (with their exponents incremented by one: 15 and 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at row 16 - this shifts the old row 16
# (INITIAL_STATE_REQ) and row 17 (FINAL_STATE_REQ) down to rows 17 and 18,
# and Excel automatically re-points range refs (G2's MAX(C2:C17) becomes
# MAX(C2:C18), the C/D/E shared-formula ranges grow to include the new
# row, etc).
$ws.Rows("16").Insert()

# Populate the new row with the SpikeList entry. Its exponent (column B)
# keeps the value that used to belong to row 16 (14) -- the rows below it
# get bumped by one instead.
$ws.Range("A16").Value2 = "SPIKE_LIST_REQ"
$ws.Range("B16").Value2 = 14
$ws.Range("C16").Formula = "=LEN(A16)"
$ws.Range("D16").Formula = '=CONCATENATE(A16, REPT(" ",$G$2+1-C16))'
$ws.Range("E16").Formula = '=CONCATENATE(D16," = (1 << ",B16,"), ")'

# The rows that got pushed down keep their names/formulas (handled by the
# Insert shift) but their exponent numbers need to move up by one.
$ws.Range("B17").Value2 = 15
$ws.Range("B18").Value2 = 16

# Scroll the view down a bit and leave the selection on the new
# SpikeList row's Final-String cell, matching where the edit was made.
$ws.Range("E16").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

$wb.Application.Calculate() | Out-Null
